# Weekly update: insert a new "Fruta / Palta" record as row 51, pushing the
# existing rows 51-69 down to 52-70 (dimension grows from A1:T69 to A1:T70).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 51 (shifts 51..69 -> 52..70).
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new record.
$ws.Range("A51").Value = 1
$ws.Range("B51").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C51").Value = "Arica y Parinacota"
$ws.Range("D51").Value = 44524
$ws.Range("E51").Value = 15
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100106
$ws.Range("H51").Value = "Oleaginosos"
$ws.Range("I51").Value = 100106002
$ws.Range("J51").Value = "Palta"
$ws.Range("K51").Value = "Hass"
$ws.Range("L51").Value = "Segunda"
$ws.Range("M51").Value = 200
$ws.Range("N51").Value = 68000
$ws.Range("O51").Value = 70000
$ws.Range("P51").Value = 69000
$ws.Range("Q51").Value = "$/caja 25 kilos"
$ws.Range("R51").Value = "Región de Coquimbo"
$ws.Range("S51").Value = 2760
$ws.Range("T51").Value = 25
